$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the K11:K19 formulas: change "=" comparison to "-" subtraction
# so the formula computes MOD(J-H*7, $I$7) instead of MOD(J=H*7, $I$7)
$ws.Range("K11").Formula = "=MOD(J11-H11*7,`$I`$7)"
$ws.Range("K12:K19").Formula = "=MOD(J12-H12*7,`$I`$7)"

# Update the selection on the sheet view to I22
$ws.Range("I22").Select()
